$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.096.78"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.16%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.637.01"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.91%  "

# Row 4
$ws.Range("E4").Value = "  -0.16%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.59%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5257"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.07%  "

# Row 7
$ws.Range("E7").Value = "  -0.14%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2600"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.32%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06319"
$ws.Range("D9").Style = "Normal"

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.73"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.17%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07663"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.28%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.639.70"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.65%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.424"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.53%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.860.94"
$ws.Range("D14").Style = "Normal"

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5502"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.12%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8176"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.01%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.09"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.72%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.078.20"
$ws.Range("D18").Style = "Normal"

# Row 19
$ws.Range("E19").Value = "  -0.12%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.686"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.42%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "188.29"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.80%  "

# Row 22
$ws.Range("E22").Value = "  -2.26%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.165"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.17%  "

# Row 24
$ws.Range("E24").Value = "  -0.10%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.81"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.61%  "

# Row 26
$ws.Range("E26").Value = "  -2.66%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.412"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.28%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.83"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.87%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.407"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.00%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05995"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.16%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.257"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.86%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.442"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.25%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.408"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.40%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.641"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.55%  "

# Row 35
$ws.Range("E35").Value = "  -1.20%  "

# Row 36
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.395"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.65%  "

# Row 37
$ws.Range("B37").Value = "MXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.762"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.16%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5740"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.26%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01617"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.04%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8556"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.05%  "

# Row 41
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.738"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.63%  "

# Row 42
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.001"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.24%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.037.23"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.84%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.62"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.34%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.787.24"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.86%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₈107"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.27%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.73"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.45%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.002"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.33%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.070"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.04%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05171"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.27%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4221"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.62%  "

Write-Host "Applied cryptos update"
